$d = $word.ActiveDocument

# Update the title/date paragraph (first paragraph, before the table)
$d.Paragraphs.Item(1).Range.Text = "2025-08-01 Friday"

# Update each table cell value (20 rows x 5 columns), in row-major order
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "91-76=15"
$t.Cell(1,2).Range.Text = "49-15=34"
$t.Cell(1,3).Range.Text = "85-74=11"
$t.Cell(1,4).Range.Text = "46-3=43"
$t.Cell(1,5).Range.Text = "13+22=35"

$t.Cell(2,1).Range.Text = "18+40=58"
$t.Cell(2,2).Range.Text = "59-4=55"
$t.Cell(2,3).Range.Text = "6+29=35"
$t.Cell(2,4).Range.Text = "47-38=9"
$t.Cell(2,5).Range.Text = "87-30=57"

$t.Cell(3,1).Range.Text = "89-67=22"
$t.Cell(3,2).Range.Text = "43+32=75"
$t.Cell(3,3).Range.Text = "6+85=91"
$t.Cell(3,4).Range.Text = "98-98=0"
$t.Cell(3,5).Range.Text = "88-54=34"

$t.Cell(4,1).Range.Text = "61-48=13"
$t.Cell(4,2).Range.Text = "71-33=38"
$t.Cell(4,3).Range.Text = "15+77=92"
$t.Cell(4,4).Range.Text = "18+20=38"
$t.Cell(4,5).Range.Text = "29+26=55"

$t.Cell(5,1).Range.Text = "83-38=45"
$t.Cell(5,2).Range.Text = "85-21=64"
$t.Cell(5,3).Range.Text = "50+42=92"
$t.Cell(5,4).Range.Text = "28+25=53"
$t.Cell(5,5).Range.Text = "33+37=70"

$t.Cell(6,1).Range.Text = "23-15=8"
$t.Cell(6,2).Range.Text = "78+8=86"
$t.Cell(6,3).Range.Text = "45-18=27"
$t.Cell(6,4).Range.Text = "55-29=26"
$t.Cell(6,5).Range.Text = "70+5=75"

$t.Cell(7,1).Range.Text = "97-28=69"
$t.Cell(7,2).Range.Text = "46-17=29"
$t.Cell(7,3).Range.Text = "24+67=91"
$t.Cell(7,4).Range.Text = "52+35=87"
$t.Cell(7,5).Range.Text = "99-27=72"

$t.Cell(8,1).Range.Text = "47-23=24"
$t.Cell(8,2).Range.Text = "6-6=0"
$t.Cell(8,3).Range.Text = "32-28=4"
$t.Cell(8,4).Range.Text = "20+15=35"
$t.Cell(8,5).Range.Text = "69-28=41"

$t.Cell(9,1).Range.Text = "6+64=70"
$t.Cell(9,2).Range.Text = "83-39=44"
$t.Cell(9,3).Range.Text = "6+13=19"
$t.Cell(9,4).Range.Text = "63-60=3"
$t.Cell(9,5).Range.Text = "64-30=34"

$t.Cell(10,1).Range.Text = "47-33=14"
$t.Cell(10,2).Range.Text = "7+52=59"
$t.Cell(10,3).Range.Text = "33-20=13"
$t.Cell(10,4).Range.Text = "26-10=16"
$t.Cell(10,5).Range.Text = "2+86=88"

$t.Cell(11,1).Range.Text = "68-30=38"
$t.Cell(11,2).Range.Text = "61+29=90"
$t.Cell(11,3).Range.Text = "27-2=25"
$t.Cell(11,4).Range.Text = "54+11=65"
$t.Cell(11,5).Range.Text = "75-48=27"

$t.Cell(12,1).Range.Text = "57-26=31"
$t.Cell(12,2).Range.Text = "14-3=11"
$t.Cell(12,3).Range.Text = "32-8=24"
$t.Cell(12,4).Range.Text = "65-60=5"
$t.Cell(12,5).Range.Text = "10-8=2"

$t.Cell(13,1).Range.Text = "16-0=16"
$t.Cell(13,2).Range.Text = "83-82=1"
$t.Cell(13,3).Range.Text = "96-92=4"
$t.Cell(13,4).Range.Text = "69-52=17"
$t.Cell(13,5).Range.Text = "45+20=65"

$t.Cell(14,1).Range.Text = "95-19=76"
$t.Cell(14,2).Range.Text = "91-64=27"
$t.Cell(14,3).Range.Text = "5+56=61"
$t.Cell(14,4).Range.Text = "87+5=92"
$t.Cell(14,5).Range.Text = "85-72=13"

$t.Cell(15,1).Range.Text = "54-32=22"
$t.Cell(15,2).Range.Text = "26+61=87"
$t.Cell(15,3).Range.Text = "7+63=70"
$t.Cell(15,4).Range.Text = "49-20=29"
$t.Cell(15,5).Range.Text = "14+82=96"

$t.Cell(16,1).Range.Text = "2+0=2"
$t.Cell(16,2).Range.Text = "91-76=15"
$t.Cell(16,3).Range.Text = "39-20=19"
$t.Cell(16,4).Range.Text = "59-58=1"
$t.Cell(16,5).Range.Text = "91-52=39"

$t.Cell(17,1).Range.Text = "9+59=68"
$t.Cell(17,2).Range.Text = "89-57=32"
$t.Cell(17,3).Range.Text = "35+44=79"
$t.Cell(17,4).Range.Text = "82+4=86"
$t.Cell(17,5).Range.Text = "20+12=32"

$t.Cell(18,1).Range.Text = "96-28=68"
$t.Cell(18,2).Range.Text = "45-30=15"
$t.Cell(18,3).Range.Text = "68+10=78"
$t.Cell(18,4).Range.Text = "17+20=37"
$t.Cell(18,5).Range.Text = "21-19=2"

$t.Cell(19,1).Range.Text = "94-89=5"
$t.Cell(19,2).Range.Text = "54-34=20"
$t.Cell(19,3).Range.Text = "11+36=47"
$t.Cell(19,4).Range.Text = "83-35=48"
$t.Cell(19,5).Range.Text = "33+0=33"

$t.Cell(20,1).Range.Text = "68-30=38"
$t.Cell(20,2).Range.Text = "4+5=9"
$t.Cell(20,3).Range.Text = "93-90=3"
$t.Cell(20,4).Range.Text = "89-86=3"
$t.Cell(20,5).Range.Text = "8+26=34"

